# "Generate Report for Handback"
#
# This handback run:
#   - marks both locales' Status as handed-back/in-sync
#   - stamps zh-cn's Latest Handback DateTime
#   - stamps de-de's Latest Handback DateTime (a later timestamp)
#   - fills in the (previously empty) "Latest Target File" column for both
#     locales, turning it into a hyperlink back to the source markdown file,
#     and records the locale-specific xliff file name in "Latest Handback File"
#   - widens a few columns on all three sheets to fit the new/longer text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Column width helper.
# This runtime's ColumnWidth setter rounds the input to 2 decimal places
# before adding the fixed 5/6 character padding that ends up in the saved
# <col width="..."/>. Back the input off by that padding so the stored
# width lands as close as possible to the desired value.
# ---------------------------------------------------------------------
function Set-ColWidth($ws, $colIndex, $desiredStoredWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $desiredStoredWidth - (5 / 6)
}

# =======================================================================
# Overview sheet
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

# E2/F2 display the shared "Ready for handoff" status text -> update it.
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"

# Columns E and F need to be wider to fit the longer status text.
Set-ColWidth $ov 5 29.9777047293527
Set-ColWidth $ov 6 29.9777047293527

# =======================================================================
# zh-cn sheet
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# Status column.
$zh.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Target File (I2) - now populated + hyperlinked to the source md.
$zh.Hyperlinks.Add(
    $zh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d452e8762b5068ce030150b454b2b206b2e9247d/e2e/a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.md",
    "",
    "",
    "a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.md"
) | Out-Null

# Latest Handback File (J2).
$zh.Range("J2").Value = "a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.277b1a8ef1fdd8287878bb058a9a454b8ec1dad1.zh-cn.xlf"

# Latest Handback DateTime (K2).
$zh.Range("K2").Value = "2016-08-25 17:01:45"

# Widen columns to fit the new content.
Set-ColWidth $zh 3 29.9777047293527
Set-ColWidth $zh 9 40
Set-ColWidth $zh 10 40

# =======================================================================
# de-de sheet
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

# Status column.
$de.Range("C2").Value = "Handed back: in sync with en-US"

# Latest Target File (I2) - now populated + hyperlinked to the source md.
$de.Hyperlinks.Add(
    $de.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d452e8762b5068ce030150b454b2b206b2e9247d/e2e/a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.md",
    "",
    "",
    "a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.md"
) | Out-Null

# Latest Handback File (J2).
$de.Range("J2").Value = "a36eee9f-1c81-4cf0-b9a7-7ed205148f7a.277b1a8ef1fdd8287878bb058a9a454b8ec1dad1.de-de.xlf"

# Latest Handback DateTime (K2).
$de.Range("K2").Value = "2016-08-25 17:01:54"

# Widen columns to fit the new content.
Set-ColWidth $de 3 29.9777047293527
Set-ColWidth $de 9 40
Set-ColWidth $de 10 40

Write-Host "Handback report generated."
